$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new columns before D, shifting existing D:K data to F:M ---
$ws.Columns("D:E").Insert()

# --- Step 2: copy number formatting from (now-shifted) F:G into the new blank D:E columns ---
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$cols = @("D","E","F","G","H","I","J","K","L","M")

# --- Step 3: "Period Ending" date rows (quarter-end dates, newest two added at D/E) ---
$dateVals = @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
foreach ($r in @(7, 38, 80)) {
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Range($cols[$i] + "$r").Value = $dateVals[$i]
    }
}

# --- Step 4: financial data rows ---
$rows = @{
    8 = @(317800,284700,246100,298400,258400,234100,216500,265700,232500,164200)
    9 = @(212300,202100,176000,204300,174000,162900,159500,181500,158700,115500)
    10 = @(105500,82500,70100,94100,84400,71100,57000,84200,73800,48600)
    12 = @("NA","NA","NA","NA","NA","NA","NA","NA","NA","NA")
    13 = @(0,0,0,0,0,0,0,0,0,0)
    14 = @(149400,0,900,"NA",18100,0,0,0,7900,-3400)
    15 = @(73700,67500,66000,65000,70100,67400,71300,74600,79000,74500)
    17 = @(440800,291800,266700,290900,279000,246900,245700,275500,263600,202900)
    18 = @(-123000,-7200,-20600,7600,-20600,-12900,-29200,-9900,-31000,-38800)
    20 = @(-2300,800,0,-800,500,700,700,800,1200,1500)
    21 = @(-51500,61100,45400,71700,50100,55200,42800,65500,49100,37300)
    22 = @(24100,23400,24300,23700,27800,24100,25800,26200,33000,27100)
    23 = @(-149300,-29700,-44900,-16900,-47900,-36400,-54300,-35200,-62800,-64300)
    24 = @(-1700,-6900,-9700,-3500,-12900,-16800,-27400,-17900,-39300,-27900)
    25 = @(0,0,0,0,0,0,0,0,0,0)
    26 = @(-147600,-22800,-35100,-13500,-35000,-19600,-26900,-17400,-23500,-36400)
    27 = @(-147600,-22800,-35100,-13500,-35000,-19600,-26900,-17400,-23500,-36400)
    28 = @(0,0,0,0,0,0,0,0,0,0)
    29 = @(0,0,0,0,0,0,0,0,0,0)
    30 = @(0,0,0,0,0,0,0,0,0,0)
    31 = @(0,0,0,0,0,0,0,0,0,0)
    32 = @(2300,-800,0,800,-500,-700,-700,-800,-1200,-1500)
    33 = @(-147600,-22800,-35100,-13500,-35000,-19600,-26900,-17400,-23500,-36400)
    34 = @(0,0,0,0,0,0,0,0,0,0)
    35 = @(-147600,-22800,-35100,-13500,-35000,-19600,-26900,-17400,-23500,-36400)
    41 = @(71900,81700,70500,60900,48400,98000,70700,92600,88900,270600)
    42 = @(0,0,0,0,0,0,0,0,0,0)
    43 = @(277100,254700,235100,286000,262000,244000,242200,275100,254900,178500)
    44 = @(25400,23900,25100,19900,18300,18800,19200,19600,18500,19500)
    45 = @(14600,"NA","NA","NA","NA","NA","NA",0,0,0)
    46 = @(389000,360200,330700,366800,328800,360900,332100,387400,362300,468600)
    47 = @(1800,1700,1800,1700,1700,"NA","NA","NA",0,0)
    48 = @(2261400,2251000,2329900,2345300,2362000,2437800,2547300,2724800,2797800,2788800)
    49 = @(26300,179600,180700,179500,173600,154400,155600,161600,161900,161500)
    50 = @(0,0,0,0,0,0,0,0,0,0)
    51 = @(0,0,0,0,0,0,0,0,0,0)
    52 = @(27400,24900,28300,31200,31100,"NA","NA","NA","NA","NA")
    53 = @(0,0,0,0,0,0,0,0,0,0)
    54 = @(2706000,2817500,2871300,2924500,2897200,2953000,3035000,3273800,3322000,3418900)
    57 = @(204300,191600,184700,165800,156000,162800,150500,196200,184900,167300)
    58 = @(0,0,0,0,0,0,0,0,0,0)
    59 = @(5700,2700,"NA","NA","NA","NA","NA",0,0,0)
    60 = @(210000,194300,184700,165800,156000,162800,150500,196200,184900,167300)
    61 = @(1269800,1264200,1291800,1322300,1287800,1323000,1372900,1454100,1465000,1540600)
    62 = @(66900,71700,80400,95100,106100,84900,102000,135400,164700,193200)
    63 = @(0,0,0,0,0,0,0,0,0,0)
    64 = @(0,0,0,0,0,0,0,0,0,0)
    65 = @(0,0,0,0,0,0,0,0,0,0)
    66 = @(1546700,1530100,1557000,1583200,1549900,1570700,1625400,1785600,1814600,1901100)
    68 = @(0,0,0,0,0,0,0,0,0,0)
    69 = @(0,0,0,0,0,0,0,0,0,0)
    70 = @(0,0,0,0,0,0,0,0,0,0)
    71 = @(0,0,0,0,0,0,0,0,0,0)
    72 = @(-728500,-580900,-558100,-522900,-509500,-474500,-454900,-441900,-424500,-401000)
    73 = @(0,0,0,0,0,0,0,0,0,0)
    74 = @(0,0,0,0,0,0,0,0,0,0)
    75 = @(0,0,0,0,0,0,0,0,0,0)
    76 = @(1159300,1287400,1314300,1341400,1347300,1382300,1409500,1488200,1507400,1517800)
    77 = @(0,0,0,0,0,0,0,0,0,0)
    81 = @(-147600,-22800,-35100,-13500,-35000,-19600,-26900,-17400,-23500,-36400)
    83 = @(73700,67500,66000,65000,70100,67400,71300,74600,79000,74500)
    84 = @(0,0,0,0,0,0,0,0,0,0)
    85 = @(0,0,0,0,0,0,0,0,0,0)
    86 = @(0,0,0,0,0,0,0,0,0,0)
    87 = @(0,0,0,0,0,0,0,0,0,0)
    88 = @(0,0,0,0,0,0,0,0,0,0)
    89 = @(69600,23800,96500,28400,17300,42200,2000,25900,-21400,13500)
    91 = @(-22000,-21100,-25500,-16600,-1300,-23400,-14500,-17000,-34400,-60000)
    92 = @(0,0,0,0,0,0,0,0,0,0)
    93 = @(0,0,0,0,0,0,0,0,0,0)
    94 = @(-14500,-11300,-31500,-17800,-13000,-13700,-20100,-21700,-41500,-51200)
    96 = @(0,0,0,0,0,0,0,0,0,0)
    97 = @(0,0,0,0,0,0,0,0,0,0)
    98 = @(0,0,0,0,0,0,0,0,0,0)
    99 = @(0,0,0,0,0,0,0,0,0,0)
    100 = @(-69000,200,-57000,0,-54700,0,0,-300,-118900,-43000)
    101 = @(4100,-1500,1600,1800,700,-1300,-900,-200,100,1200)
    102 = @(-9800,11200,9500,12500,-49600,27300,-19000,3700,-181700,-79500)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Range($cols[$i] + "$r").Value = $vals[$i]
    }
}